$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (new) = Row 4 (old) ---
$ws.Range("A2").Value2 = 91961108
$ws.Range("B2").Value2 = 78569
$ws.Range("E2").Value2 = 6458
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value2 = "Lunglav"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value2 = "Lobaria pulmonaria"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value2 = "(L.) Hoffm."
$ws.Range("H2").Style = "Normal"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value2 = "Mysjöberget, Mpd"
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Value2 = 632736.7829380766
$ws.Range("R2").Value2 = 6940262.09546657
$ws.Range("S2").Value2 = 25
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = "2020-06-09"
$ws.Range("Y2").Style = "Normal"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = "2020-06-09"
$ws.Range("AA2").Style = "Normal"
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value2 = "Mikael Gudrunsson"
$ws.Range("AW2").Style = "Normal"
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value2 = "Mikael Gudrunsson"
$ws.Range("AX2").Style = "Normal"
$ws.Range("AR2").ClearContents()

# --- Row 3 (new) = Row 2 (old) ---
$ws.Range("A3").Value2 = 66541020
$ws.Range("B3").Value2 = 89392
$ws.Range("E3").Value2 = 1202
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "Ullticka"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value2 = ""
$ws.Range("I3").Style = "Normal"
$ws.Range("AR3").NumberFormat = "@"
$ws.Range("AR3").Value2 = "23349"
$ws.Range("AR3").Style = "Normal"

# --- Row 4 (new) = Row 3 (old) ---
$ws.Range("A4").Value2 = 66541021
$ws.Range("B4").Value2 = 77506
$ws.Range("E4").Value2 = 6425
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value2 = "Garnlav"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value2 = "Alectoria sarmentosa"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value2 = "(Ach.) Ach."
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value2 = "10"
$ws.Range("I4").Style = "Normal"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value2 = "Grysjöbäcken, Mpd"
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value2 = 633480.2438334802
$ws.Range("R4").Value2 = 6940243.949297423
$ws.Range("S4").Value2 = 50
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2014-09-18"
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2014-09-18"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AR4").NumberFormat = "@"
$ws.Range("AR4").Value2 = "23350"
$ws.Range("AR4").Style = "Normal"
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value2 = "Malin Sahlin"
$ws.Range("AW4").Style = "Normal"
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value2 = "Via Malin Sahlin"
$ws.Range("AX4").Style = "Normal"

